$d = $word.ActiveDocument

# Locate the table row describing option "5" (note-range trigger, filtered
# by velocity) in the first table of the document, then highlight the
# three data cells on that row ("5 Trigger note range...", "Note Max",
# "Note Min") in yellow - matching Word's "Text Highlight Color" tool.

$tbl = $d.Tables.Item(1)
$targetRow = -1
for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    $cellText = $tbl.Cell($r, 2).Range.Text
    if ($cellText -like "*5 Trigger note range*") {
        $targetRow = $r
        break
    }
}

if ($targetRow -gt 0) {
    foreach ($col in 2, 3, 4) {
        $cell = $tbl.Cell($targetRow, $col)
        $cell.Range.Font.HighlightColorIndex = 7
    }
}
